$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two empty placeholder cells (AA275, AD275) from row 275
$ws.Range("AA275").ClearContents()
$ws.Range("AD275").ClearContents()

# Row 276
$ws.Cells.Item(276, 1).Value = "'2023-01-26 09:27:33"
$ws.Cells.Item(276, 2).Value = "'18.9"
$ws.Cells.Item(276, 3).Value = "'26.6"
$ws.Cells.Item(276, 4).Value = "'4.1"
$ws.Cells.Item(276, 5).Value = "'."
$ws.Cells.Item(276, 6).Value = "'24.6"
$ws.Cells.Item(276, 7).Value = "'."
$ws.Cells.Item(276, 8).Value = "'41"
$ws.Cells.Item(276, 9).Value = "'."
$ws.Cells.Item(276, 10).Value = "'50"
$ws.Cells.Item(276, 11).Value = "'50"
$ws.Cells.Item(276, 12).Value = "'0"
$ws.Cells.Item(276, 13).Value = "'0"
$ws.Cells.Item(276, 14).Value = "'50"
$ws.Cells.Item(276, 15).Value = "'50"
$ws.Cells.Item(276, 16).Value = "'."
$ws.Cells.Item(276, 17).Value = "'."
$ws.Cells.Item(276, 18).Value = "'100"
$ws.Cells.Item(276, 19).Value = "'0"
$ws.Cells.Item(276, 20).Value = "'0"
$ws.Cells.Item(276, 21).Value = "'0"
$ws.Cells.Item(276, 22).Value = "'100"
$ws.Cells.Item(276, 23).Value = "'0"
$ws.Cells.Item(276, 24).Value = "'100"
$ws.Cells.Item(276, 25).Value = "'."
$ws.Cells.Item(276, 26).Value = "'9.22"
$ws.Cells.Item(276, 28).Value = "'7002200"
$ws.Cells.Item(276, 29).Value = "'2200"
$ws.Cells.Item(276, 31).Value = "'19.6"
$ws.Cells.Item(276, 32).Value = "'26.5"
$ws.Cells.Item(276, 33).Value = "'4.1"
$ws.Cells.Item(276, 34).Value = "'50"
$ws.Cells.Item(276, 35).Value = "'50"
$ws.Cells.Item(276, 36).Value = "'50"
$ws.Cells.Item(276, 37).Value = "'50"
$ws.Cells.Item(276, 38).Value = "'50"
$ws.Cells.Item(276, 39).Value = "'0"
$ws.Cells.Item(276, 40).Value = "'."
$ws.Cells.Item(276, 41).Value = "'619"
$ws.Cells.Item(276, 42).Value = "'32"
$ws.Cells.Item(276, 43).Value = "'0"
$ws.Cells.Item(276, 44).Value = "'0.20"
$ws.Cells.Item(276, 45).Value = "'0.17"
$ws.Cells.Item(276, 46).Value = "'0.15"
$ws.Cells.Item(276, 47).Value = "'162.02"
$ws.Cells.Item(276, 48).Value = "'4341.50"
$ws.Cells.Item(276, 49).Value = "'0.00"
$ws.Cells.Item(276, 50).Value = "'144.46"
$ws.Cells.Item(276, 51).Value = "'3560.09"
$ws.Cells.Item(276, 52).Value = "'0.53"
$ws.Cells.Item(276, 53).Value = "'413.74"
$ws.Cells.Item(276, 54).Value = "'12292.90"
$ws.Cells.Item(276, 55).Value = "'19.0"
$ws.Cells.Item(276, 56).Value = "'."
$ws.Cells.Item(276, 57).Value = "'."
$ws.Cells.Item(276, 58).Value = "'41"
$ws.Cells.Item(276, 59).Value = "'."
$ws.Cells.Item(276, 60).Value = "'41"
$ws.Cells.Item(276, 61).Value = "'146030596"

# Row 277
$ws.Cells.Item(277, 1).Value = "'2023-01-26 11:33:29"
$ws.Cells.Item(277, 2).Value = "'18.8"
$ws.Cells.Item(277, 3).Value = "'26.5"
$ws.Cells.Item(277, 4).Value = "'4.2"
$ws.Cells.Item(277, 5).Value = "'."
$ws.Cells.Item(277, 6).Value = "'24.4"
$ws.Cells.Item(277, 7).Value = "'."
$ws.Cells.Item(277, 8).Value = "'41"
$ws.Cells.Item(277, 9).Value = "'."
$ws.Cells.Item(277, 10).Value = "'50"
$ws.Cells.Item(277, 11).Value = "'50"
$ws.Cells.Item(277, 12).Value = "'0"
$ws.Cells.Item(277, 13).Value = "'0"
$ws.Cells.Item(277, 14).Value = "'50"
$ws.Cells.Item(277, 15).Value = "'50"
$ws.Cells.Item(277, 16).Value = "'."
$ws.Cells.Item(277, 17).Value = "'."
$ws.Cells.Item(277, 18).Value = "'100"
$ws.Cells.Item(277, 19).Value = "'0"
$ws.Cells.Item(277, 20).Value = "'0"
$ws.Cells.Item(277, 21).Value = "'0"
$ws.Cells.Item(277, 22).Value = "'100"
$ws.Cells.Item(277, 23).Value = "'1"
$ws.Cells.Item(277, 24).Value = "'100"
$ws.Cells.Item(277, 25).Value = "'."
$ws.Cells.Item(277, 26).Value = "'9.11"
$ws.Cells.Item(277, 28).Value = "'7002200"
$ws.Cells.Item(277, 29).Value = "'2200"
$ws.Cells.Item(277, 31).Value = "'19.8"
$ws.Cells.Item(277, 32).Value = "'26.5"
$ws.Cells.Item(277, 33).Value = "'4.2"
$ws.Cells.Item(277, 34).Value = "'50"
$ws.Cells.Item(277, 35).Value = "'50"
$ws.Cells.Item(277, 36).Value = "'50"
$ws.Cells.Item(277, 37).Value = "'50"
$ws.Cells.Item(277, 38).Value = "'50"
$ws.Cells.Item(277, 39).Value = "'1"
$ws.Cells.Item(277, 40).Value = "'."
$ws.Cells.Item(277, 41).Value = "'0"
$ws.Cells.Item(277, 42).Value = "'30"
$ws.Cells.Item(277, 43).Value = "'0"
$ws.Cells.Item(277, 44).Value = "'0.21"
$ws.Cells.Item(277, 45).Value = "'0.17"
$ws.Cells.Item(277, 46).Value = "'0.22"
$ws.Cells.Item(277, 47).Value = "'162.09"
$ws.Cells.Item(277, 48).Value = "'4341.57"
$ws.Cells.Item(277, 49).Value = "'0.00"
$ws.Cells.Item(277, 50).Value = "'144.46"
$ws.Cells.Item(277, 51).Value = "'3560.09"
$ws.Cells.Item(277, 52).Value = "'0.94"
$ws.Cells.Item(277, 53).Value = "'414.15"
$ws.Cells.Item(277, 54).Value = "'12293.31"
$ws.Cells.Item(277, 55).Value = "'19.0"
$ws.Cells.Item(277, 56).Value = "'."
$ws.Cells.Item(277, 57).Value = "'."
$ws.Cells.Item(277, 58).Value = "'41"
$ws.Cells.Item(277, 59).Value = "'."
$ws.Cells.Item(277, 60).Value = "'41"
$ws.Cells.Item(277, 61).Value = "'144982020"

# Row 278
$ws.Cells.Item(278, 1).Value = "'2023-01-30 18:22:48"
$ws.Cells.Item(278, 2).Value = "'18.3"
$ws.Cells.Item(278, 3).Value = "'27.2"
$ws.Cells.Item(278, 4).Value = "'4.0"
$ws.Cells.Item(278, 5).Value = "'."
$ws.Cells.Item(278, 6).Value = "'24.9"
$ws.Cells.Item(278, 7).Value = "'."
$ws.Cells.Item(278, 8).Value = "'45"
$ws.Cells.Item(278, 9).Value = "'."
$ws.Cells.Item(278, 10).Value = "'20"
$ws.Cells.Item(278, 11).Value = "'20"
$ws.Cells.Item(278, 12).Value = "'0"
$ws.Cells.Item(278, 13).Value = "'0"
$ws.Cells.Item(278, 14).Value = "'20"
$ws.Cells.Item(278, 15).Value = "'20"
$ws.Cells.Item(278, 16).Value = "'."
$ws.Cells.Item(278, 17).Value = "'."
$ws.Cells.Item(278, 18).Value = "'0"
$ws.Cells.Item(278, 19).Value = "'0"
$ws.Cells.Item(278, 20).Value = "'0"
$ws.Cells.Item(278, 21).Value = "'0"
$ws.Cells.Item(278, 22).Value = "'100"
$ws.Cells.Item(278, 23).Value = "'3"
$ws.Cells.Item(278, 24).Value = "'."
$ws.Cells.Item(278, 25).Value = "'."
$ws.Cells.Item(278, 26).Value = "'10.30"
$ws.Cells.Item(278, 28).Value = "'7002200"
$ws.Cells.Item(278, 29).Value = "'2200"
$ws.Cells.Item(278, 31).Value = "'18.2"
$ws.Cells.Item(278, 32).Value = "'27.2"
$ws.Cells.Item(278, 33).Value = "'4.0"
$ws.Cells.Item(278, 34).Value = "'20"
$ws.Cells.Item(278, 35).Value = "'20"
$ws.Cells.Item(278, 36).Value = "'20"
$ws.Cells.Item(278, 37).Value = "'20"
$ws.Cells.Item(278, 38).Value = "'20"
$ws.Cells.Item(278, 39).Value = "'3"
$ws.Cells.Item(278, 40).Value = "'."
$ws.Cells.Item(278, 41).Value = "'0"
$ws.Cells.Item(278, 42).Value = "'10"
$ws.Cells.Item(278, 43).Value = "'0"
$ws.Cells.Item(278, 44).Value = "'0.15"
$ws.Cells.Item(278, 45).Value = "'0.16"
$ws.Cells.Item(278, 46).Value = "'0.20"
$ws.Cells.Item(278, 47).Value = "'131.43"
$ws.Cells.Item(278, 48).Value = "'4343.77"
$ws.Cells.Item(278, 49).Value = "'0.00"
$ws.Cells.Item(278, 50).Value = "'114.43"
$ws.Cells.Item(278, 51).Value = "'3560.09"
$ws.Cells.Item(278, 52).Value = "'0.31"
$ws.Cells.Item(278, 53).Value = "'359.54"
$ws.Cells.Item(278, 54).Value = "'12309.62"
$ws.Cells.Item(278, 55).Value = "'19.0"
$ws.Cells.Item(278, 56).Value = "'."
$ws.Cells.Item(278, 57).Value = "'."
$ws.Cells.Item(278, 58).Value = "'45"
$ws.Cells.Item(278, 59).Value = "'."
$ws.Cells.Item(278, 60).Value = "'45"
$ws.Cells.Item(278, 61).Value = "'137502724"

# Row 279
$ws.Cells.Item(279, 1).Value = "'2023-01-30 18:29:48"
$ws.Cells.Item(279, 2).Value = "'18.5"
$ws.Cells.Item(279, 3).Value = "'27.3"
$ws.Cells.Item(279, 4).Value = "'4.0"
$ws.Cells.Item(279, 5).Value = "'."
$ws.Cells.Item(279, 6).Value = "'24.9"
$ws.Cells.Item(279, 7).Value = "'."
$ws.Cells.Item(279, 8).Value = "'45"
$ws.Cells.Item(279, 9).Value = "'."
$ws.Cells.Item(279, 10).Value = "'20"
$ws.Cells.Item(279, 11).Value = "'20"
$ws.Cells.Item(279, 12).Value = "'0"
$ws.Cells.Item(279, 13).Value = "'0"
$ws.Cells.Item(279, 14).Value = "'20"
$ws.Cells.Item(279, 15).Value = "'20"
$ws.Cells.Item(279, 16).Value = "'."
$ws.Cells.Item(279, 17).Value = "'."
$ws.Cells.Item(279, 18).Value = "'0"
$ws.Cells.Item(279, 19).Value = "'0"
$ws.Cells.Item(279, 20).Value = "'0"
$ws.Cells.Item(279, 21).Value = "'0"
$ws.Cells.Item(279, 22).Value = "'100"
$ws.Cells.Item(279, 23).Value = "'3"
$ws.Cells.Item(279, 24).Value = "'."
$ws.Cells.Item(279, 25).Value = "'."
$ws.Cells.Item(279, 26).Value = "'10.30"
$ws.Cells.Item(279, 27).Value = "x"
$ws.Cells.Item(279, 27).Value = ""
$ws.Cells.Item(279, 28).Value = "'7002200"
$ws.Cells.Item(279, 29).Value = "'2200"
$ws.Cells.Item(279, 30).Value = "x"
$ws.Cells.Item(279, 30).Value = ""
$ws.Cells.Item(279, 31).Value = "'18.6"
$ws.Cells.Item(279, 32).Value = "'27.2"
$ws.Cells.Item(279, 33).Value = "'3.9"
$ws.Cells.Item(279, 34).Value = "'20"
$ws.Cells.Item(279, 35).Value = "'20"
$ws.Cells.Item(279, 36).Value = "'20"
$ws.Cells.Item(279, 37).Value = "'20"
$ws.Cells.Item(279, 38).Value = "'20"
$ws.Cells.Item(279, 39).Value = "'3"
$ws.Cells.Item(279, 40).Value = "'."
$ws.Cells.Item(279, 41).Value = "'0"
$ws.Cells.Item(279, 42).Value = "'10"
$ws.Cells.Item(279, 43).Value = "'0"
$ws.Cells.Item(279, 44).Value = "'0.15"
$ws.Cells.Item(279, 45).Value = "'0.16"
$ws.Cells.Item(279, 46).Value = "'0.20"
$ws.Cells.Item(279, 47).Value = "'131.44"
$ws.Cells.Item(279, 48).Value = "'4343.77"
$ws.Cells.Item(279, 49).Value = "'0.00"
$ws.Cells.Item(279, 50).Value = "'114.43"
$ws.Cells.Item(279, 51).Value = "'3560.09"
$ws.Cells.Item(279, 52).Value = "'0.32"
$ws.Cells.Item(279, 53).Value = "'359.55"
$ws.Cells.Item(279, 54).Value = "'12309.62"
$ws.Cells.Item(279, 55).Value = "'19.0"
$ws.Cells.Item(279, 56).Value = "'."
$ws.Cells.Item(279, 57).Value = "'."
$ws.Cells.Item(279, 58).Value = "'45"
$ws.Cells.Item(279, 59).Value = "'."
$ws.Cells.Item(279, 60).Value = "'45"
$ws.Cells.Item(279, 61).Value = "'137502724"
